$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells to remain text (matching original inline-string cells)
# so Excel does not auto-convert numeric-looking values (e.g. "355.12") into numbers.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '51.652.53'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '2.799.86'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '355.12'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = '109.53'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').Value = '0.558'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  +5.60%  '
$ws.Range('D10').Value = '40.15'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').Value = '0.0839'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '20.00'
$ws.Range('E13').Value = '  +2.95%  '
$ws.Range('D14').Value = '7.79'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '3.238.58'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = '2.793.28'
$ws.Range('E16').Value = '  +0.98%  '
$ws.Range('D17').Value = '0.941'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '51.639.72'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = '7.76'
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('D20').Value = '3.17'
$ws.Range('E20').Value = '  +2.87%  '
$ws.Range('D21').Value = '13.37'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('D24').Value = '268.22'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '2.78'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '26.10'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').Value = '0.165'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('D30').Value = '37.56'
$ws.Range('E30').Value = '  +9.27%  '
$ws.Range('E31').Value = '  +4.89%  '
$ws.Range('D32').Value = '6.36'
$ws.Range('E32').Value = '  +10.96%  '
$ws.Range('D33').Value = '52.21'
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('D34').Value = '5.66'
$ws.Range('E34').Value = '  +9.14%  '
$ws.Range('E35').Value = '  -4.92%  '
$ws.Range('D36').Value = '0.0856'
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '18.81'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('D39').Value = '3.15'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('D40').Value = '1.99'
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('E42').Value = '  -4.98%  '
$ws.Range('D43').Value = '119.85'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('D45').Value = '21.87'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '3.45'
$ws.Range('E46').Value = '  +5.80%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.140.53'
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('D48').Value = '2.37'
$ws.Range('E48').Value = '  +7.07%  '
$ws.Range('D49').Value = '0.923'
$ws.Range('E49').Value = '  -3.56%  '
$ws.Range('D50').Value = '1.37'
$ws.Range('E50').Value = '  +11.02%  '
$ws.Range('D51').Value = '0.222'
$ws.Range('E51').Value = '  +16.90%  '
